$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-23 Sunday" "2025-03-24 Monday"

Replace-Text "34×34=1156" "41×15=615"
Replace-Text "95×68=6460" "52×62=3224"
Replace-Text "98×64=6272" "87×84=7308"
Replace-Text "97×24=2328" "89×13=1157"
Replace-Text "60×43=2580" "21×72=1512"

Replace-Text "66×74=4884" "93×91=8463"
Replace-Text "59×25=1475" "94×11=1034"
Replace-Text "86×24=2064" "44×21=924"
Replace-Text "75×83=6225" "45×80=3600"
Replace-Text "68×68=4624" "30×72=2160"

Replace-Text "45×32=1440" "26×88=2288"
Replace-Text "99×94=9306" "21×74=1554"
Replace-Text "47×13=611" "52×25=1300"
Replace-Text "22×17=374" "34×14=476"
Replace-Text "49×95=4655" "13×75=975"

Replace-Text "70×99=6930" "94×23=2162"
Replace-Text "97×58=5626" "95×36=3420"
Replace-Text "39×12=468" "18×25=450"
Replace-Text "38×62=2356" "38×55=2090"
Replace-Text "85×36=3060" "35×66=2310"

Replace-Text "65×81=5265" "92×71=6532"
Replace-Text "65×44=2860" "53×73=3869"
Replace-Text "61×24=1464" "39×96=3744"
Replace-Text "60×55=3300" "30×53=1590"
Replace-Text "41×56=2296" "57×62=3534"
